$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $c = $ws.Range($cellRef)
    $c.Formula = "'" + $value
    $c.Style = "Normal"
}

Set-TextValue 'D2' '301.60'
Set-TextValue 'E2' '1.87%'
Set-TextValue 'D3' '43.95'
Set-TextValue 'E3' '6.80%'
Set-TextValue 'D4' '5.073'
Set-TextValue 'E4' '0.71%'
Set-TextValue 'D5' '0.07695'
Set-TextValue 'E5' '3.57%'
$ws.Range('B6').Value = 'FTXToken'
$ws.Range('C6').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextValue 'D6' '1.616'
Set-TextValue 'E6' '2.87%'
$ws.Range('B7').Value = 'MXToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue 'D7' '1.041'
Set-TextValue 'E7' '12.68%'
$ws.Range('B8').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C8').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue 'D8' '0.1277'
Set-TextValue 'E8' '9.08%'
$ws.Range('B9').Value = 'WazirX'
$ws.Range('C9').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue 'D9' '0.1870'
Set-TextValue 'E9' '2.93%'
$ws.Range('B10').Value = 'MandalaExchangeToken'
$ws.Range('C10').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue 'D10' '0.09187'
Set-TextValue 'E10' '4.24%'
$ws.Range('B11').Value = 'BitrueCoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue 'D11' '0.04170'
Set-TextValue 'E11' '-2.95%'
$ws.Range('B12').Value = 'BitMartToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue 'D12' '0.1049'
Set-TextValue 'E12' '-0.20%'
$ws.Range('B13').Value = 'BitForexToken'
$ws.Range('C13').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue 'D13' '0.001284'
Set-TextValue 'E13' '0.40%'
$ws.Range('B14').Value = 'TigerCash'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue 'D14' '0.005757'
Set-TextValue 'E14' '-2.32%'
$ws.Range('B15').Value = 'UpBots'
$ws.Range('C15').Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
Set-TextValue 'D15' '0.007430'
Set-TextValue 'E15' '1,894.35%'
$ws.Range('B16').Value = 'LEO'
$ws.Range('C16').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue 'D16' '3.343'
Set-TextValue 'E16' '-0.38%'
$ws.Range('B17').Value = 'GateToken'
$ws.Range('C17').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextValue 'D17' '4.417'
Set-TextValue 'E17' '1.38%'
Set-TextValue 'D18' '2.330'
Set-TextValue 'E18' '-3.70%'
Set-TextValue 'D19' '0.3352'
Set-TextValue 'E19' '1.99%'
Set-TextValue 'D20' '8.671'
Set-TextValue 'E20' '9.92%'
Set-TextValue 'D21' '0.1399'
Set-TextValue 'E21' '2.71%'
Set-TextValue 'D22' '0.3175'
Set-TextValue 'E22' '6.92%'
Set-TextValue 'D23' '0.04192'
Set-TextValue 'E23' '3.88%'
Set-TextValue 'D24' '0.001285'
Set-TextValue 'E24' '0.97%'
Set-TextValue 'D25' '0.004473'
Set-TextValue 'E25' '15.74%'
Set-TextValue 'D26' '0.0001349'
Set-TextValue 'D38' '0.02496'
Set-TextValue 'E38' '4.66%'
Set-TextValue 'D39' '0.05297'
Set-TextValue 'E39' '2.05%'
Set-TextValue 'D40' '0.005932'
Set-TextValue 'E40' '-10.24%'
Set-TextValue 'D41' '0.007681'
Set-TextValue 'E41' '-1.28%'
Set-TextValue 'D42' '0.1348'
Set-TextValue 'E42' '2.49%'
Set-TextValue 'D43' '0.007362'
Set-TextValue 'E43' '-0.35%'
Set-TextValue 'D44' '0.007556'
Set-TextValue 'E44' '-3.13%'
Set-TextValue 'D45' '0.3000'
Set-TextValue 'E45' '-6.62%'
Set-TextValue 'D46' '0.00006653'
Set-TextValue 'E46' '6.44%'
Set-TextValue 'E47' '-0.18%'
Set-TextValue 'D48' '0.04155'
Set-TextValue 'E48' '-10.18%'
Set-TextValue 'E49' '-0.08%'
Set-TextValue 'E50' '-0.18%'
Set-TextValue 'E51' '-0.18%'
